$d = $word.ActiveDocument

$replacements = @(
    @{old="380×2="; new="314×9="},
    @{old="779×4="; new="324×9="},
    @{old="477×3="; new="817×5="},
    @{old="224×2="; new="151×2="},
    @{old="922×2="; new="838×7="},
    @{old="584×6="; new="626×6="},
    @{old="136×7="; new="936×7="},
    @{old="490×2="; new="522×3="},
    @{old="295×2="; new="740×7="},
    @{old="962×3="; new="379×9="},
    @{old="623×8="; new="274×2="},
    @{old="895×3="; new="451×8="},
    @{old="356×5="; new="290×7="},
    @{old="123×6="; new="983×5="},
    @{old="519×3="; new="530×2="},
    @{old="768×5="; new="181×9="},
    @{old="989×4="; new="280×3="},
    @{old="843×7="; new="219×4="},
    @{old="512×7="; new="415×5="},
    @{old="595×3="; new="553×4="},
    @{old="571×8="; new="307×5="},
    @{old="797×7="; new="471×2="},
    @{old="321×9="; new="779×8="},
    @{old="402×6="; new="271×5="},
    @{old="222×6="; new="578×5="}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $r.new, 2)
}
